$wb = $excel.ActiveWorkbook
$wsR = $wb.Worksheets.Item("REZISTIVITE")
$wsS = $wb.Worksheets.Item("SISMIK")
$wsD = $wb.Worksheets.Item("SONDAJ")

# --- REZISTIVITE (sheet1): shift the Z-depth values (E,G,I,K) left by one
# layer slot; the R-value columns (F,H,J,L) stay put. Last slot becomes blank.
$rez = @{
    2 = @{ E = 17; G = 40; I = $null }
    3 = @{ E = 9;  G = 43; I = $null }
    4 = @{ E = 11; G = 24; I = 46;  K = $null }
    5 = @{ E = 12; G = 29; I = 50;  K = $null }
    6 = @{ E = 10; G = 27; I = 46;  K = $null }
    7 = @{ E = 13; G = 30; I = $null }
}
foreach ($r in $rez.Keys) {
    foreach ($col in $rez[$r].Keys) {
        $val = $rez[$r][$col]
        $addr = "$col$r"
        if ($null -eq $val) {
            $wsR.Range($addr).ClearContents()
        } else {
            $wsR.Range($addr).Value = $val
        }
    }
}

# --- SISMIK (sheet2): same left-shift for the Z-depth values (E,H,K); the
# 4th layer columns (Z4/VP4/VS4 = N,O,P) are removed entirely.
$sis = @{
    2 = @{ E = 18; H = 30; K = $null }
    3 = @{ E = 18; H = 39; K = $null }
    4 = @{ E = 13; H = 30; K = $null }
    5 = @{ E = 11; H = 25; K = $null }
    6 = @{ E = 11; H = 28; K = $null }
    7 = @{ E = 4;  H = 19; K = $null }
    8 = @{ E = 13; H = $null }
    9 = @{ E = 15; H = $null }
}
foreach ($r in $sis.Keys) {
    foreach ($col in $sis[$r].Keys) {
        $val = $sis[$r][$col]
        $addr = "$col$r"
        if ($null -eq $val) {
            $wsS.Range($addr).ClearContents()
        } else {
            $wsS.Range($addr).Value = $val
        }
    }
}

# Remove the now-unused 4th-layer columns (Z4(m), VP4(m/s), VS4(m/s)) entirely.
$wsS.Range("N1:P1").EntireColumn.Delete()

# --- View state: SISMIK becomes the active tab, with L9 selected on both
# REZISTIVITE and SISMIK (SONDAJ's selection is unchanged).
[void]$wsR.Range("L9").Select()
[void]$wsS.Range("L9").Select()
$wsS.Activate()
